$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header: O1 "Medicinsk" ---
$ws.Range("O1").Value = 'Medicinsk'

# --- 10 new plant rows (72-81) ---
# Row 72: Surkörsbär (Prunus cerasus)
$ws.Range("A72").Value = 'Prunus cerasus'
$ws.Range("B72").Value = 'Surkörsbär'
$ws.Range("C72").Value = 'Surkörsbär (Prunus cerasus) är en art i familjen rosväxter. Arten förekommer förvildad i största delen av Europa. Anses fungera bättre till sylt, saft och vin eftersom bäret har starkare arom och hållbarhet än odlade varianter. Går att äta råa men kan vara ganska sura.'
$ws.Range("D72").Value = 2
$ws.Range("J72").Value = 1
$ws.Range("L72").Value = 1
$ws.Range("N72").Value = 1

# Row 73: Sötkörsbär (Prunus avium)
$ws.Range("A73").Value = 'Prunus avium'
$ws.Range("B73").Value = 'Sötkörsbär'
$ws.Range("C73").Value = 'Sötkörsbär eller fågelbär, Prunus avium, är en art i plommonsläktet inom familjen rosväxter. Arten förekommer naturligt i Europa. Det är träd som blir cirka 15–30 meter högt, och den art från vilken de flesta odlade kultivarer av körsbär tagits fram. Det finns i three main variteter såsom Bigarråer, vanlig sötkörsbär'
$ws.Range("A73").Font.Italic = $true
$ws.Range("D73").Value = 2
$ws.Range("J73").Value = 1
$ws.Range("N73").Value = 1

# Row 74: Hägg (Prunus padus)
$ws.Range("A74").Value = 'Prunus padus'
$ws.Range("B74").Value = 'Hägg'
$ws.Range("C74").Value = 'Barken och frukterna är rika på amygdalin vilket bittermandelolja kan framställas ur. Bären är sträva ungefär som slånbär men är ätliga. Med fördel kokas saft eller vin på bäret. Kärnorna är svåra att skilja från bäret och används därför sällan till sylt.'
$ws.Range("A74").Font.Italic = $true
$ws.Range("D74").Value = 2
$ws.Range("J74").Value = 1
$ws.Range("N74").Value = 1

# Row 75: Slån (Prunus spinosa)
$ws.Range("A75").Value = 'Prunus spinosa'
$ws.Range("B75").Value = 'Slån'
$ws.Range("C75").Value = 'Den får blå eller blåsvarta stenfrukter som kallas slånbär. Busken, som har vassa tornar, växer vilt i större delen av Europa. I Sverige är den vanlig i södra Sveriges kustlandskap. Bären innehåller garvsyra som ger dem en kärv och sur, nästan bitter smak. När temperaturen sjunker under noll grader omvandlas garvsyran, och slånbärens smak blir sötare och mindre sträv. Det fungerar utmärkt att imitera detta genom att lägga dem i frysen under 1–2 dygn. Bären har länge använts till att koka saft och vin. Bladen och blommorna kan användas till te. Veden är mycket hårt och har använts till sniderier. Barken till garvning av läder.'
$ws.Range("A75").Font.Italic = $true
$ws.Range("D75").Value = 2
$ws.Range("G75").Value = 2
$ws.Range("H75").Value = 2
$ws.Range("J75").Value = 1
$ws.Range("N75").Value = 2

# Row 76: Brudbröd (Filipendula vulgaris)
$ws.Range("A76").Value = 'Filipendula vulgaris'
$ws.Range("B76").Value = 'Brudbröd'
$ws.Range("C76").Value = 'Brudbröd är lätt igenkänd på de starka och sega rötternas tjocka, kulformiga uppsvällningar, hårda men dock något köttiga knölar, som magasinerar vatten för torra perioder. Brudbrödet är Älgörtens närmaste släkting. Rötterna kan kokas och ätas som potatis. '
$ws.Range("A76").Font.Italic = $true
$ws.Range("D76").Value = 2
$ws.Range("E76").Value = 2

# Row 77: Älggräs (Filipendula ulmaria)
$ws.Range("A77").Value = 'Filipendula ulmaria'
$ws.Range("B77").Value = 'Älggräs'
$ws.Range("C77").Value = 'Älggräs (Filipendula ulmaria L.), älgört, är en ört i släktet älggräs i familjen rosväxter. Den högväxta arten har små vita och starkt doftande blommor. Blommor och blad innehåller acetylsalicylsyra, och arten har använts både i medicinskt syfte och som smaksättning av drycker. Den växer i fuktiga marker såsom diken och skuggiga stränder. Kan finnas i mängder på ängar. Blommor kan användas som te och fungerar mot huvudvärk då den innehåller salicylsyra. '
$ws.Range("D77").Value = 1
$ws.Range("H77").Value = 2
$ws.Range("O77").Value = 'Huvudvärk'

# Row 78: Måbär (Ribes alpinum)
$ws.Range("A78").Value = 'Ribes alpinum'
$ws.Range("B78").Value = 'Måbär'
$ws.Range("C78").Value = 'Måbärsbusken är ganska allmän i lundar och skogsbackar, i synnerhet i östra Götaland och Svealand, men finns glest utbredd från Skåne till mellersta Norrland. Han och honblommor sitter på skilda buskar. Liknar vinbär och bären är ätbara, söta men lite fadda i smaken.'
$ws.Range("A78").Font.Italic = $true
$ws.Range("D78").Value = 2
$ws.Range("J78").Value = 1
$ws.Range("N78").Value = 2

# Row 79: Röda vinbär (Ribes rubrum)
$ws.Range("A79").Value = 'Ribes rubrum'
$ws.Range("B79").Value = 'Röda vinbär'
$ws.Range("C79").Value = 'Röda vinbär är ett samlingsnamn för trädgårdsvinbär (Ribes rubrum) och skogsvinbär (Ribes spicatum). Ibland har bägge dessa ansetts vara underarter av samma huvudart. Förvildade trädgårdsvinbär kan förekomma i samma områden som skogsvinbär. Vita vinbär är en variant. Bären är sura och kan ätas råa eller kokas till sylt, saft eller gelé. Bären innehåller stora mängder pektin och kan därför användas tillsammans med andra bär för att skapa gelé.'
$ws.Range("D79").Value = 2
$ws.Range("J79").Value = 1
$ws.Range("L79").Value = 1
$ws.Range("N79").Value = 2

# Row 80: Svarta vinbär (Ribes nigrum)
$ws.Range("A80").Value = 'Ribes nigrum'
$ws.Range("B80").Value = 'Svarta vinbär'
$ws.Range("C80").Value = 'Svarta vinbär eller svartvinbär (Ribes nigrum) är en växt som tillhör Vinbärssläktet och familjen ripsväxter. En gammal synonym är tistron. Bären är mycket rika på C-vitamin och anses nyttigare och godare än det röda vinbäret. Används till sylt, saft, vin eller gelé. Bladen kan användas som fläderblommor eller till te och gurkinläggningar.'
$ws.Range("D80").Value = 2
$ws.Range("G80").Value = 2
$ws.Range("J80").Value = 1
$ws.Range("L80").Value = 1
$ws.Range("N80").Value = 2

# Row 81: Krusbär (Ribes uva-crispa)
$ws.Range("A81").Value = 'Ribes uva-crispa'
$ws.Range("B81").Value = 'Krusbär'
$ws.Range("C81").Value = 'Krusbär (Ribes uva-crispa) är namnet både på en buske och dess bär. Den dök förmodligen upp i sverige på 1500-talet. Busken har taggiga grenar. Bären kan drabbas av krusbärsmjöldaggen som är en parasitsvamp som täcker in bäret i en brunt ludd och gör den oätlig. Bären används till sylt, saft, gelé eller vin.'
$ws.Range("C81").Characters(10, 16).Font.Italic = $true
$ws.Range("D81").Value = 2
$ws.Range("J81").Value = 1
$ws.Range("L81").Value = 1
$ws.Range("N81").Value = 2

# --- Column O width (new column) ---
$ws.Columns.Item(15).ColumnWidth = 10.6

# --- Sheet view: scroll frozen pane and move selection ---
$ws.Activate()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.ScrollRow = 59
$ws.Range("A79").Select()
